$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "_old" / "_new" header suffixes to "_FV2410" / "_FV2504" ---
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- Turn the data range into a real table (ListObject) ---
# Stash a copy of row 1's formatting off to the side first so we can
# restore it exactly after ListObjects.Add (which otherwise recolors the
# header row to the table style and bakes the previous look into a dxf).
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A70:U70")
$headerRange.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$headerRange.ClearFormats() | Out-Null
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), 0, 1)

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$scratch.Clear() | Out-Null

$lo.TableStyle = $null

# --- Freeze the header row ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
